$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.714.84"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "2.270.60"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "503.31"
$ws.Range("D6").Value = "127.68"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("E8").Value = "  +0.46%  "
$ws.Range("D9").Value = "2.282.50"
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("D10").Value = "0.0974"
$ws.Range("E10").Value = "  +2.41%  "
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").Value = "5.08"
$ws.Range("E12").Value = "  +7.57%  "
$ws.Range("D13").Value = "0.339"
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("D14").Value = "23.40"
$ws.Range("E14").Value = "  +3.33%  "
$ws.Range("D15").Value = "2.676.19"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "54.819.85"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "2.274.08"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").Value = "10.34"
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("D21").Value = "309.81"
$ws.Range("E21").Value = "  +2.30%  "
$ws.Range("D22").Value = "6.56"
$ws.Range("E22").Value = "  +3.81%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "59.78"
$ws.Range("E24").Value = "  -2.10%  "
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("E26").Value = "  +3.11%  "
$ws.Range("E27").Value = "  +2.26%  "
$ws.Range("D28").Value = "172.14"
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("D29").Value = "6.10"
$ws.Range("E29").Value = "  +3.11%  "
$ws.Range("E30").Value = "  +1.21%  "
$ws.Range("D31").Value = "0.0₃0700"
$ws.Range("E31").Value = "  +1.37%  "
$ws.Range("E32").Value = "  +5.38%  "
$ws.Range("D34").Value = "17.90"
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("D35").Value = "0.996"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("E36").Value = "  +2.27%  "
$ws.Range("D37").Value = "0.900"
$ws.Range("E37").Value = "  -4.50%  "
$ws.Range("E38").Value = "  +4.21%  "
$ws.Range("D39").Value = "36.73"
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("E40").Value = "  +2.98%  "
$ws.Range("D41").Value = "0.374"
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("D42").Value = "134.68"
$ws.Range("E42").Value = "  +7.90%  "
$ws.Range("D43").Value = "3.45"
$ws.Range("E43").Value = "  +3.00%  "
$ws.Range("D44").Value = "4.83"
$ws.Range("E44").Value = "  +0.55%  "
$ws.Range("D45").Value = "256.88"
$ws.Range("E45").Value = "  +7.55%  "
$ws.Range("E46").Value = "  +2.35%  "
$ws.Range("E47").Value = "  +2.37%  "
$ws.Range("D48").Value = "0.544"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").Value = "0.373"
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("D51").Value = "16.36"
$ws.Range("E51").Value = "  +1.37%  "
